# SU_Suspension/params.xlsx - "Design table pour inserts"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value: Rayon_rotule (mm) goes from 10 to 8
$ws.Range("B10").Value = 8

# New rows for the insert design table
$ws.Range("A11").Value = "Longueur_collage_plug (mm)"
$ws.Range("B11").Value = 35

$ws.Range("A12").Value = "Epaisseur_rotule (mm)"
$ws.Range("B12").Value = 16

$ws.Range("A13").Value = "Plug_rint (mm)"
$ws.Range("B13").Value = 6.5

$ws.Range("A14").Value = "Plug_rext (mm)"
$ws.Range("B14").Formula = "=15.7/2"

# Match new selection state
$ws.Range("B12").Select()
